# Update countries & provincias Spain
# Daily data refresh: update the "last updated" timestamp, refresh several
# countries' COVID-19 stats, and re-sort a couple of rows whose totals moved
# past their neighbours (Suecia overtakes Arabia Saudita/Irlanda; Mauritania
# overtakes Butan).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header / timestamp note (row 1)
$ws.Range("A1").Value = "Datos actualizados a 29 de Abril de 2020 a las 14:22"

# Row 17 - Paises Bajos (stats refresh only)
$ws.Range("A17").Value = "Paises Bajos"
$ws.Range("B17").Value = 38802
$ws.Range("C17").Value = 386
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 33841
$ws.Range("F17").Value = 861
$ws.Range("G17").Value = 145
$ws.Range("H17").Value = 4711

# Rows 23-25 - Suecia overtakes Arabia Saudita and Irlanda
$ws.Range("A23").Value = "Suecia"
$ws.Range("B23").Value = 20302
$ws.Range("C23").Value = 681
$ws.Range("D23").Value = 1005
$ws.Range("E23").Value = 16835
$ws.Range("F23").Value = 479
$ws.Range("G23").Value = 107
$ws.Range("H23").Value = 2462

$ws.Range("A24").Value = "Arabia Saudita"
$ws.Range("B24").Value = 20077
$ws.Range("C24").Value = 0
$ws.Range("D24").Value = 2784
$ws.Range("E24").Value = 17141
$ws.Range("F24").Value = 118
$ws.Range("G24").Value = 0
$ws.Range("H24").Value = 152

$ws.Range("A25").Value = "Irlanda"
$ws.Range("B25").Value = 19877
$ws.Range("C25").Value = 0
$ws.Range("D25").Value = 9233
$ws.Range("E25").Value = 9485
$ws.Range("F25").Value = 141
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 1159

# Row 30 - Pakistan (stats refresh only)
$ws.Range("A30").Value = "Pakistan"
$ws.Range("B30").Value = 15289
$ws.Range("C30").Value = 677
$ws.Range("D30").Value = 3425
$ws.Range("E30").Value = 11529
$ws.Range("F30").Value = 111
$ws.Range("G30").Value = 23
$ws.Range("H30").Value = 335

# Row 41 - Dinamarca (stats refresh only)
$ws.Range("A41").Value = "Dinamarca"
$ws.Range("B41").Value = 9008
$ws.Range("C41").Value = 157
$ws.Range("D41").Value = 6366
$ws.Range("E41").Value = 2199
$ws.Range("F41").Value = 66
$ws.Range("G41").Value = 9
$ws.Range("H41").Value = 443

# Row 61 - Kazajistan (stats refresh only)
$ws.Range("A61").Value = "Kazajistan"
$ws.Range("B61").Value = 3105
$ws.Range("C61").Value = 78
$ws.Range("D61").Value = 798
$ws.Range("E61").Value = 2282
$ws.Range("F61").Value = 41
$ws.Range("G61").Value = 0
$ws.Range("H61").Value = 25

# Row 67 - Croacia (stats refresh only)
$ws.Range("A67").Value = "Croacia"
$ws.Range("B67").Value = 2062
$ws.Range("C67").Value = 15
$ws.Range("D67").Value = 1288
$ws.Range("E67").Value = 707
$ws.Range("F67").Value = 19
$ws.Range("G67").Value = 4
$ws.Range("H67").Value = 67

# Row 95 - Albania (stats refresh only)
$ws.Range("A95").Value = "Albania"
$ws.Range("B95").Value = 766
$ws.Range("C95").Value = 16
$ws.Range("D95").Value = 455
$ws.Range("E95").Value = 281
$ws.Range("F95").Value = 4
$ws.Range("G95").Value = 0
$ws.Range("H95").Value = 30

# Rows 208-209 - Mauritania overtakes Butan
$ws.Range("A208").Value = "Mauritania"
$ws.Range("B208").Value = 8
$ws.Range("C208").Value = 1
$ws.Range("D208").Value = 6
$ws.Range("E208").Value = 1
$ws.Range("F208").Value = 0
$ws.Range("G208").Value = 0
$ws.Range("H208").Value = 1

$ws.Range("A209").Value = "Butan"
$ws.Range("B209").Value = 7
$ws.Range("C209").Value = 0
$ws.Range("D209").Value = 5
$ws.Range("E209").Value = 2
$ws.Range("F209").Value = 0
$ws.Range("G209").Value = 0
$ws.Range("H209").Value = 0
